# Generate Report for Handoff
#
# Row 3 of every sheet tracks the file "b.md". A new handoff has been
# produced for it, so:
#   - the Status column moves from "Handed back: in sync with en-US"
#     to "Ready for handoff" (Overview + zh-cn + de-de)
#   - the per-locale sheets get a new "Latest Handoff File" name and a
#     new "Latest Handoff Datetime" for that row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: row 3 ("b.md") status for both locales
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------
# zh-cn sheet: row 3 ("b.md") status + latest handoff file/datetime
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-03-10 18:31:00"

# The hyperlink on C3 keeps pointing at the same handoff-history page,
# only the displayed file name changes. Rebuild the sheet's hyperlinks
# collection (in place edits aren't exposed) so the C3 link shows the
# new file name while every other link is preserved unchanged.
$zhcn.Range("A1").Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a886c338e329e523e7e8bfb2764cc0568521d73d/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a73df203666331c6eb1e37ec7017d5ee8c6e0c11/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/bb789de700bc611a50ccee485de2651db6932f4d/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b84e3b309982e09db0e31685e622f0ff4e85af39/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a886c338e329e523e7e8bfb2764cc0568521d73d/e2e/b.md", "", "", "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a73df203666331c6eb1e37ec7017d5ee8c6e0c11/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/bb789de700bc611a50ccee485de2651db6932f4d/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b84e3b309982e09db0e31685e622f0ff4e85af39/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a886c338e329e523e7e8bfb2764cc0568521d73d/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------
# de-de sheet: row 3 ("b.md") status + latest handoff file/datetime
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-03-10 18:31:06"

$dede.Range("A1").Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a886c338e329e523e7e8bfb2764cc0568521d73d/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/406978126a66925839daa0c239d601c30a7d4515/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ae7ac16c94db6afcbc7afec37610bdffbf1c87cb/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3d6c50458cbd5acd9c5975f8ea78562a58499cda/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a886c338e329e523e7e8bfb2764cc0568521d73d/e2e/b.md", "", "", "b.md")
$dede.Hyperlinks.Add($dede.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/406978126a66925839daa0c239d601c30a7d4515/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ae7ac16c94db6afcbc7afec37610bdffbf1c87cb/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3d6c50458cbd5acd9c5975f8ea78562a58499cda/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a886c338e329e523e7e8bfb2764cc0568521d73d/.localization-config", "", "", ".localization-config")
